$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.614.45'
$ws.Range("E2").Value = '  +4.88%  '
$ws.Range("D3").Value = '3.633.58'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.01'
$ws.Range("E5").Value = '  +1.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '193.33'
$ws.Range("E6").Value = '  +4.59%  '
$ws.Range("E7").Value = '  +2.14%  '
$ws.Range("D8").Value = '3.627.71'
$ws.Range("E8").Value = '  +4.79%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.179'
$ws.Range("E10").Value = '  +3.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.671'
$ws.Range("E11").Value = '  +3.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '58.48'
$ws.Range("E12").Value = '  +3.66%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000291'
$ws.Range("E13").Value = '  +4.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.92'
$ws.Range("E14").Value = '  +4.88%  '
$ws.Range("D15").Value = '4.210.70'
$ws.Range("E15").Value = '  +5.00%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.77'
$ws.Range("E16").Value = '  +5.83%  '
$ws.Range("D17").Value = '3.627.70'
$ws.Range("E17").Value = '  +4.97%  '
$ws.Range("D18").Value = '70.556.05'
$ws.Range("E18").Value = '  +4.91%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.71'
$ws.Range("E19").Value = '  +4.97%  '
$ws.Range("E20").Value = '  +0.43%  '
$ws.Range("E21").Value = '  +4.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '488.95'
$ws.Range("E22").Value = '  +1.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.40'
$ws.Range("E23").Value = '  +13.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.41'
$ws.Range("E24").Value = '  -1.96%  '
$ws.Range("E25").Value = '  +0.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '91.07'
$ws.Range("E26").Value = '  +0.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.15'
$ws.Range("E27").Value = '  +6.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.38'
$ws.Range("E28").Value = '  +3.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.70'
$ws.Range("E29").Value = '  +5.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.11'
$ws.Range("E30").Value = '  +5.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.85'
$ws.Range("E31").Value = '  +9.40%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.120'
$ws.Range("E32").Value = '  +7.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.31'
$ws.Range("E33").Value = '  +4.81%  '
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '625.78'
$ws.Range("E34").Value = '  +5.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '66.03'
$ws.Range("E35").Value = '  +2.67%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '39.68'
$ws.Range("E36").Value = '  +8.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.415'
$ws.Range("E37").Value = '  +7.40%  '
$ws.Range("D38").Value = '0.0₃0819'
$ws.Range("E38").Value = '  +6.12%  '
$ws.Range("E39").Value = '  -1.44%  '
$ws.Range("E41").Value = '  +1.19%  '
$ws.Range("D42").Value = '3.300.72'
$ws.Range("E42").Value = '  +2.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.17'
$ws.Range("E43").Value = '  +8.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.79'
$ws.Range("E44").Value = '  +9.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0452'
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.139'
$ws.Range("E46").Value = '  +2.71%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.28'
$ws.Range("E47").Value = '  +2.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.24'
$ws.Range("E48").Value = '  +5.66%  '
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.34'
$ws.Range("E50").Value = '  +3.11%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '143.37'
$ws.Range("E51").Value = '  +2.11%  '
